# Fix truncated text in 2026-02-11 detailed analysis slide (Sources list).
# Restore full sentence-boundary source citations in the "Sources" textbox
# (slide 11, shape "TextBox 3"). Each citation paragraph is rewritten in
# place via TextRange.Paragraphs(Start, Length) sub-ranges, which preserves
# every paragraph/run formatting property (pPr/defRPr) untouched.
#
# The text is written in two steps (placeholder, then final value): the
# host's text-replacement diffing can otherwise keep a trailing/leading
# character that is identical between the old and new string as a leftover
# un-replaced run (e.g. two URLs that both end in "/"), splitting a single
# citation into extra <a:r> runs. Routing every edit through a placeholder
# that shares no characters with either string sidesteps that.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(3)

$sources = @(
    'NBC News: https://www.nbcnews.com/tech/internet/openai-starts-testing-ads-chatgpt-rcna258242',
    'Intellizence: https://intellizence.com/insights/startup-funding/startup-funding-trends-january-2026-ai-infrastructure-and-robotics/',
    'Tech Xplore: https://techxplore.com/news/2026-02-ai-limits-generative-video.html',
    'National Law Review: https://natlawreview.com/article/2026-outlook-artificial-intelligence',
    'Medium: https://medium.com/@urano10/the-future-of-ai-models-in-2026-whats-actually-coming-410141f3c979',
    'DigiTimes: https://www.digitimes.com/news/a20251215PD230/meta-ai-llm-llama-development.html',
    'TechCrunch: https://techcrunch.com/2026/02/10/ai-video-startup-runway-raises-315m-at-5-3b-valuation-eyes-more-capable-world-models/',
    'Axios: https://www.axios.com/2026/02/10/ai-ceo-feuds-openai-anthropic-google',
    'MIT Technology Review: https://www.technologyreview.com/2026/01/05/1130662/whats-next-for-ai-in-2026/',
    'TechCrunch: https://techcrunch.com/2026/01/02/in-2026-ai-will-move-from-hype-to-pragmatism/',
    'Digital Watch Observatory: https://dig.watch/updates/adobe-firefly-unlocks-powerful-unlimited-ai-generation-in-2026',
    'InfoQ: https://www.infoq.com/news/2026/01/microsoft-llm-contextual-privacy/'
)

for ($i = 0; $i -lt $sources.Count; $i++) {
    # Paragraph 1 is the blank spacer paragraph before the citation list,
    # so citation $i (0-based) lives at paragraph index $i + 2.
    $paraIndex = $i + 2

    $tr = $shp.TextFrame.TextRange
    $placeholder = $tr.Paragraphs($paraIndex, 1)
    $placeholder.Text = "."

    $tr = $shp.TextFrame.TextRange
    $para = $tr.Paragraphs($paraIndex, 1)
    $para.Text = $sources[$i]
}

